$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Header-row borders ------------------------------------------------
# C1/D1 (and, on sheet 2, F1/G1) currently inherit style s=1 (bold + full
# thin border + centered alignment) from the merged header cell (B1:D1 /
# E1:G1). The target instead gives them a plain/unstyled cell that only
# carries a top+bottom border (interior column) or a top+bottom+right
# border (rightmost column, closing the header box).
#
# Build the two canonical styles once on sheet 1, then fan them out with
# a formats-only paste so every cell that needs the same look lands on
# the exact same shared cellXfs entry instead of minting duplicates.
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Anonymize "fedcore" column headers --------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty placeholder cell -----------------------------
$ws2.Range("G5").ClearContents()
